$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 2550
$ws.Range("I40").Value = 2550
$ws.Range("K40").Value = 2550
$ws.Range("M40").Value = -2375

# Row 119
$ws.Range("H119").Value = 993.3333
$ws.Range("J119").Value = 993.3333
$ws.Range("L119").Value = 2979.9999
$ws.Range("N119").Value = -12655.9999

# Row 120
$ws.Range("H120").Value = 72326.664
$ws.Range("J120").Value = 72326.664
$ws.Range("L120").Value = 72326.664
$ws.Range("N120").Value = -82002.664

# Row 121
$ws.Range("H121").Value = 500
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 500
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 1500
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -4994

# Row 131
$ws.Range("H131").Value = 54828.473
$ws.Range("I131").Value = 60443
$ws.Range("J131").Value = 7105
$ws.Range("K131").Value = 181329
$ws.Range("L131").Value = 21315
$ws.Range("M131").Value = -176289
$ws.Range("N131").Value = -31395


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 110
$ws.Range("H110").Value = 1924.5
$ws.Range("I110").Value = 811
$ws.Range("J110").Value = 2147.2
$ws.Range("K110").Value = 811
$ws.Range("L110").Value = 2147.2
$ws.Range("M110").Value = 1234
$ws.Range("N110").Value = -6237.2

# Row 122
$ws.Range("H122").Value = 1892.5
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 1828
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 5484
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -10384


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 59
$ws.Range("H59").Value = 37813.344
$ws.Range("J59").Value = 37813.344
$ws.Range("L59").Value = 37813.344
$ws.Range("N59").Value = -39507.344


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 107
$ws.Range("H107").Value = 654.2826
$ws.Range("I107").Value = 454.4138
$ws.Range("J107").Value = 995.2353000000001
$ws.Range("K107").Value = 454.4138
$ws.Range("L107").Value = 995.2353000000001
$ws.Range("M107").Value = 1465.5862
$ws.Range("N107").Value = -4835.2353


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 35
$ws.Range("H35").Value = 3000
$ws.Range("J35").Value = 3000
$ws.Range("L35").Value = 9000
$ws.Range("N35").Value = -9576

# Row 54
$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14441

# Row 57
$ws.Range("H57").Value = 4326.6
$ws.Range("I57").Value = 999.5
$ws.Range("J57").Value = 4838.4614
$ws.Range("K57").Value = 2998.5
$ws.Range("L57").Value = 14515.3842
$ws.Range("M57").Value = -2439.5
$ws.Range("N57").Value = -15633.3842

# Row 62
$ws.Range("H62").Value = 2024.5555
$ws.Range("I62").Value = 812
$ws.Range("J62").Value = 2371
$ws.Range("K62").Value = 2436
$ws.Range("L62").Value = 7113
$ws.Range("M62").Value = -1750
$ws.Range("N62").Value = -8485

# Row 63
$ws.Range("H63").Value = 3505.6667
$ws.Range("I63").Value = 2117
$ws.Range("J63").Value = 4200
$ws.Range("K63").Value = 6351
$ws.Range("L63").Value = 12600
$ws.Range("M63").Value = -5602
$ws.Range("N63").Value = -14098

# Row 65
$ws.Range("H65").Value = 2024.5555
$ws.Range("I65").Value = 812
$ws.Range("J65").Value = 2371
$ws.Range("K65").Value = 7308
$ws.Range("L65").Value = 21339
$ws.Range("M65").Value = -3876
$ws.Range("N65").Value = -28203

# Row 66
$ws.Range("H66").Value = 3505.6667
$ws.Range("I66").Value = 2117
$ws.Range("J66").Value = 4200
$ws.Range("K66").Value = 19053
$ws.Range("L66").Value = 37800
$ws.Range("M66").Value = -15309
$ws.Range("N66").Value = -45288

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

# Row 92
$ws.Range("H92").Value = 418.7143
$ws.Range("I92").Value = 450.33334
$ws.Range("K92").Value = 1351.00002
$ws.Range("M92").Value = -103.0000199999999

# Row 93
$ws.Range("H93").Value = 3481.5151
$ws.Range("J93").Value = 3543.4375
$ws.Range("L93").Value = 10630.3125
$ws.Range("N93").Value = -14374.3125

# Row 94
$ws.Range("H94").Value = 4316.5
$ws.Range("J94").Value = 4980
$ws.Range("L94").Value = 14940
$ws.Range("N94").Value = -16292

# Row 97
$ws.Range("H97").Value = 596.25
$ws.Range("J97").Value = 750.5714
$ws.Range("L97").Value = 2251.7142
$ws.Range("N97").Value = -3243.7142

# Row 98
$ws.Range("H98").Value = 367.66666
$ws.Range("I98").Value = 301.5
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 904.5
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 593.5
$ws.Range("N98").Value = -4496

# Row 99
$ws.Range("H99").Value = 1033.6923
$ws.Range("I99").Value = 1033.6923
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3101.0769
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -855.0769
$ws.Range("N99").ClearContents()

# Row 101
$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 5000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -19868

# Row 102
$ws.Range("H102").Value = 7483.3335
$ws.Range("J102").Value = 7483.3335
$ws.Range("L102").Value = 22450.0005
$ws.Range("N102").Value = -27318.0005

# Row 103
$ws.Range("H103").Value = 312.5
$ws.Range("I103").Value = 312.5
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 937.5
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -58.5
$ws.Range("N103").ClearContents()

# Row 116
$ws.Range("H116").Value = 2372.25
$ws.Range("I116").Value = 1280.8572
$ws.Range("J116").Value = 3900.2
$ws.Range("K116").Value = 3842.5716
$ws.Range("L116").Value = 11700.6
$ws.Range("M116").Value = -400.5715999999998
$ws.Range("N116").Value = -18584.6

# Row 130
$ws.Range("H130").Value = 1757.5
$ws.Range("J130").Value = 3000
$ws.Range("L130").Value = 9000
$ws.Range("N130").Value = -19040

# Row 131
$ws.Range("H131").Value = 903.5714
$ws.Range("I131").Value = 348.625
$ws.Range("J131").Value = 1643.5
$ws.Range("K131").Value = 1045.875
$ws.Range("L131").Value = 4930.5
$ws.Range("M131").Value = 3994.125
$ws.Range("N131").Value = -15010.5

